# "error solve ifrs list" -- correct the IFRS figures that were scraped wrong
# for 대유플러스: rows 2-6 (FY2014-2018) get restated column values, and
# rows 7-9 (FY2019E-2021E) had their data columns wiped, keeping only the
# period labels in columns A:C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 5017
$ws.Range("E2").Value = 6
$ws.Range("F2").Value = 6
$ws.Range("G2").Value = -24
$ws.Range("H2").Value = -39
$ws.Range("I2").Value = -76
$ws.Range("J2").Value = 37
$ws.Range("K2").Value = 7067
$ws.Range("L2").Value = 5997
$ws.Range("M2").Value = 1069
$ws.Range("N2").Value = 981
$ws.Range("O2").Value = 89
$ws.Range("P2").Value = 441
$ws.Range("Q2").Value = 172
$ws.Range("R2").Value = -390
$ws.Range("S2").Value = 290
$ws.Range("T2").Value = 336
$ws.Range("U2").Value = -164
$ws.Range("V2").Value = 786
$ws.Range("W2").Value = 0.12
$ws.Range("X2").Value = -0.77
$ws.Range("Y2").Value = -7.53
$ws.Range("Z2").Value = -0.57
$ws.Range("AA2").Value = 560.74
$ws.Range("AB2").Value = 98.83
$ws.Range("AC2").Value = -86
$ws.Range("AD2").Value = -14.3
$ws.Range("AE2").Value = 1113
$ws.Range("AF2").Value = 1.11
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 88120526

# Row 3
$ws.Range("D3").Value = 4958
$ws.Range("E3").Value = 15
$ws.Range("F3").Value = 15
$ws.Range("G3").Value = 21
$ws.Range("H3").Value = -2
$ws.Range("I3").Value = -79
$ws.Range("J3").Value = 78
$ws.Range("K3").Value = 7845
$ws.Range("L3").Value = 6779
$ws.Range("M3").Value = 1066
$ws.Range("N3").Value = 917
$ws.Range("O3").Value = 149
$ws.Range("P3").Value = 441
$ws.Range("Q3").Value = -52
$ws.Range("R3").Value = -133
$ws.Range("S3").Value = 94
$ws.Range("T3").Value = 174
$ws.Range("U3").Value = -226
$ws.Range("V3").Value = 870
$ws.Range("W3").Value = 0.29
$ws.Range("X3").Value = -0.03
$ws.Range("Y3").Value = -8.359999999999999
$ws.Range("Z3").Value = -0.02
$ws.Range("AA3").Value = 635.95
$ws.Range("AB3").Value = 77.36
$ws.Range("AC3").Value = -90
$ws.Range("AD3").Value = -13.21
$ws.Range("AE3").Value = 1041
$ws.Range("AF3").Value = 1.14
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 88120526

# Row 4
$ws.Range("D4").Value = 4341
$ws.Range("E4").Value = 104
$ws.Range("F4").Value = 104
$ws.Range("G4").Value = 125
$ws.Range("H4").Value = 95
$ws.Range("I4").Value = -37
$ws.Range("J4").Value = 132
$ws.Range("K4").Value = 9452
$ws.Range("L4").Value = 8253
$ws.Range("M4").Value = 1199
$ws.Range("N4").Value = 910
$ws.Range("O4").Value = 289
$ws.Range("P4").Value = 441
$ws.Range("Q4").Value = 215
$ws.Range("R4").Value = 16
$ws.Range("S4").Value = -138
$ws.Range("T4").Value = 127
$ws.Range("U4").Value = 88
$ws.Range("V4").Value = 763
$ws.Range("W4").Value = 2.39
$ws.Range("X4").Value = 2.19
$ws.Range("Y4").Value = -4.07
$ws.Range("Z4").Value = 1.1
$ws.Range("AA4").Value = 688.39
$ws.Range("AB4").Value = 79
$ws.Range("AC4").Value = -42
$ws.Range("AD4").Value = -22.38
$ws.Range("AE4").Value = 1032
$ws.Range("AF4").Value = 0.92
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 88120526

# Row 5
$ws.Range("D5").Value = 4577
$ws.Range("E5").Value = 71
$ws.Range("F5").Value = 71
$ws.Range("G5").Value = 42
$ws.Range("H5").Value = -35
$ws.Range("I5").Value = -133
$ws.Range("J5").Value = 98
$ws.Range("K5").Value = 9463
$ws.Range("L5").Value = 8261
$ws.Range("M5").Value = 1202
$ws.Range("N5").Value = 761
$ws.Range("O5").Value = 441
$ws.Range("P5").Value = 441
$ws.Range("Q5").Value = -69
$ws.Range("R5").Value = -25
$ws.Range("S5").Value = -18
$ws.Range("T5").Value = 114
$ws.Range("U5").Value = -183
$ws.Range("V5").Value = 681
$ws.Range("W5").Value = 1.56
$ws.Range("X5").Value = -0.77
$ws.Range("Y5").Value = -15.9
$ws.Range("Z5").Value = -0.37
$ws.Range("AA5").Value = 687.25
$ws.Range("AB5").Value = 52.21
$ws.Range("AC5").Value = -151
$ws.Range("AD5").Value = -4.84
$ws.Range("AE5").Value = 864
$ws.Range("AF5").Value = 0.85
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 88120526

# Row 6
$ws.Range("D6").Value = 5034
$ws.Range("E6").Value = 364
$ws.Range("F6").Value = 364
$ws.Range("G6").Value = 224
$ws.Range("H6").Value = 191
$ws.Range("I6").Value = 25
$ws.Range("K6").Value = 10041
$ws.Range("L6").Value = 8635
$ws.Range("M6").Value = 1406
$ws.Range("N6").Value = 872
$ws.Range("P6").Value = 503
$ws.Range("Q6").Value = 588
$ws.Range("R6").Value = -328
$ws.Range("S6").Value = -125
$ws.Range("T6").Value = 153
$ws.Range("U6").Value = 435
$ws.Range("V6").Value = 920
$ws.Range("W6").Value = 7.22
$ws.Range("X6").Value = 3.8
$ws.Range("Y6").Value = 3.09
$ws.Range("Z6").Value = 1.96
$ws.Range("AA6").Value = 614.11
$ws.Range("AB6").Value = 60.05
$ws.Range("AC6").Value = 27
$ws.Range("AD6").Value = 31.98
$ws.Range("AE6").Value = 866
$ws.Range("AF6").Value = 0.99
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 100682460

# Rows 7-9 no longer carry any figures beyond the period label (columns A:C)
$ws.Range("D7:AJ9").ClearContents()
